# "Knobs is now a class" - rename the knobs.* labels in column A to match
# the new `knobs.` class-qualified naming scheme, and restore the view
# state (scroll position / selection) left behind by the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Label renames (column A, rows 19-23) ---------------------------------
# knobs.knob_link_ids              -> knobs.link_ids
# force_manual_knob_boundaries     -> knobs.force_manual_knob_boundaries
# knobs.knob_boundaries_min        -> knobs.boundaries_min
# knobs.knob_boundaries_max        -> knobs.boundaries_max
# isnaive_knob_boundaries          -> knobs.isnaive_boundaries
$ws.Range("A19").Value = "knobs.link_ids"
$ws.Range("A20").Value = "knobs.force_manual_knob_boundaries"
$ws.Range("A22").Value = "knobs.boundaries_min"
$ws.Range("A23").Value = "knobs.boundaries_max"
$ws.Range("A21").Value = "knobs.isnaive_boundaries"

# --- View state -------------------------------------------------------------
# Scroll the window so row 10 is at the top and select A21 (previously the
# selection was A8:A11).
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
[void]$ws.Range("A21").Select()
